$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns to match refreshed crypto data.
# Column D values are forced as text (leading apostrophe + style reset) so that
# numeric-looking prices (e.g. "159.12") are not auto-converted to numbers by Excel,
# matching the original inline-string cell content/type.

$ws.Range('D2').Value = "'" + '68.341.08'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').Value = "'" + '2.651.33'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'" + '598.33'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('D6').Value = "'" + '159.12'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.86%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.43%  '
$ws.Range('E9').Value = '  +4.02%  '
$ws.Range('E10').Value = '  -1.14%  '
$ws.Range('E11').Value = '  +0.50%  '
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('E13').Value = '  +0.22%  '
$ws.Range('E14').Value = '  +0.97%  '
$ws.Range('D15').Value = "'" + '3.132.26'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.40%  '
$ws.Range('D16').Value = "'" + '68.312.32'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').Value = "'" + '2.631.61'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.55%  '
$ws.Range('D18').Value = "'" + '11.43'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('D19').Value = "'" + '364.67'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.54%  '
$ws.Range('E20').Value = '  -1.05%  '
$ws.Range('D21').Value = "'" + '4.44'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +3.51%  '
$ws.Range('E22').Value = '  -0.67%  '
$ws.Range('E23').Value = '  -2.38%  '
$ws.Range('D24').Value = "'" + '75.14'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.00%  '
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('E26').Value = '  -3.07%  '
$ws.Range('D27').Value = "'" + '2.784.32'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.51%  '
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('D30').Value = "'" + '559.99'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.70%  '
$ws.Range('D31').Value = "'" + '8.04'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('E32').Value = '  -0.96%  '
$ws.Range('E33').Value = '  +0.16%  '
$ws.Range('E34').Value = '  -1.14%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').Value = '  +1.62%  '
$ws.Range('D37').Value = "'" + '19.86'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +2.88%  '
$ws.Range('D38').Value = "'" + '159.54'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.47%  '
$ws.Range('D39').Value = "'" + '0.371'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.64%  '
$ws.Range('E40').Value = '  -2.60%  '
$ws.Range('D41').Value = "'" + '5.37'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.59%  '
$ws.Range('D42').Value = "'" + '0.0₆0335'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +4.63%  '
$ws.Range('E43').Value = '  -0.54%  '
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('D45').Value = "'" + '158.48'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.29%  '
$ws.Range('D46').Value = "'" + '3.77'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('D47').Value = "'" + '22.27'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.27%  '
$ws.Range('E48').Value = '  -1.12%  '
$ws.Range('E49').Value = '  -0.18%  '
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('D51').Value = "'" + '0.569'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.16%  '
